# split the PlayerAction from player with script
# r.<Method>() calls are routed through r.Action.<Method>() now.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "if(type==2){r.Action.DeleteSelectCard();return true;}return false;"
$ws.Range("C5").Value = "if(type==3){r.Action.DeleteSelectCard();return true;}return false;"
$ws.Range("C6").Value = "r.Action.DeleteSelectCard();return true;"
$ws.Range("C7").Value = "r.Action.RecostSelectCard();return true;"

$ws.Range("C5").Select()
